$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 92, shifting existing rows 92-137 down to 93-138
$ws.Rows.Item(92).Insert()

# Populate the newly inserted row 92 with the new record
$ws.Cells.Item(92,1).Value  = 10
$ws.Cells.Item(92,2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(92,3).Value  = "La Araucanía"
$ws.Cells.Item(92,4).Value  = 44755
$ws.Cells.Item(92,5).Value  = 9
$ws.Cells.Item(92,6).Value  = "Fruta"
$ws.Cells.Item(92,7).Value  = 100104
$ws.Cells.Item(92,8).Value  = "Frutos de pepita"
$ws.Cells.Item(92,9).Value  = 100104001
$ws.Cells.Item(92,10).Value = "Granada"
$ws.Cells.Item(92,11).Value = "Wonderfull"
$ws.Cells.Item(92,12).Value = "Primera"
$ws.Cells.Item(92,13).Value = 25
$ws.Cells.Item(92,14).Value = 13000
$ws.Cells.Item(92,15).Value = 13000
$ws.Cells.Item(92,16).Value = 13000
$ws.Cells.Item(92,17).Value = "$/bandeja 10 kilos granel"
$ws.Cells.Item(92,18).Value = "Provincia de Limarí"
$ws.Cells.Item(92,19).Value = 1300
$ws.Cells.Item(92,20).Value = 10
